$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("students")

$ws.Cells.Item(16, 1).Value = "iepoy"
$ws.Cells.Item(16, 2).Value = "Luiz Kieth Patiag"
$ws.Cells.Item(16, 3).Value = "ef797c8118f02dfb649607dd5d3f8c7623048c9c063d532cc95c5ed7a898a64f"
$ws.Cells.Item(16, 4).Value = 0
